$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. Numeric-looking price strings are
# prefixed with a leading apostrophe so Excel stores them as literal text
# (matching the workbook's existing text-formatted Price column) instead of
# coercing them into floating point numbers.
$updates = @{
    'D2' = '''249.43'
    'D3' = '''22.82'
    'D4' = '''5.434'
    'D5' = '''0.05621'
    'D6' = '''3.425'
    'D7' = '''6.372'
    'D8' = '''0.8117'
    'D9' = '''0.8985'
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D10' = '''0.1431'
    'E10' = '9WazirXWRX'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '''0.07520'
    'E11' = '10MandalaExchangeTokenMDX'
    'B12' = 'LiechtensteinCryptoassetsExchange'
    'C12' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D12' = '''0.03103'
    'E12' = '11LiechtensteinCryptoassetsExchangeLCX'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D13' = '''0.03097'
    'E13' = '12BitrueCoinBTR'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D14' = '''0.09322'
    'E14' = '13BitMartTokenBMX'
    'B15' = 'MCDex'
    'C15' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'D15' = '''3.559'
    'E15' = '14MCDexMCB'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D16' = '''0.001607'
    'E16' = '15BitForexTokenBF'
    'B17' = 'CoinExToken'
    'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D17' = '''0.04754'
    'E17' = '16CoinExTokenCET'
    'B18' = 'One'
    'C18' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D18' = '''0.0005796'
    'E18' = '17OneONE'
    'D20' = '''0.004989'
    'D21' = '''0.001032'
    'D23' = '''3.702'
    'D24' = '''2.183'
    'E27' = '26AAXTokenAAB'
    'D28' = '''0.0003009'
    'D40' = '''0.04072'
    'B41' = 'BKEXToken'
    'C41' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D41' = '''0.1067'
    'E41' = '40BKEXTokenBKK'
    'B42' = 'CEJI'
    'C42' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'D42' = '''0.002723'
    'E42' = '41CEJICEJI'
    'B43' = 'KickToken'
    'C43' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D43' = '''0.002942'
    'E43' = '42KickTokenKICKWorstin24h'
    'D44' = '''0.007780'
    'D45' = '''0.00005579'
    'D47' = '''0.5006'
    'D48' = '''0.2397'
    'E48' = '47BOLOBOLOBestin24h'
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Host "Applied $($updates.Count) cell updates"
